$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Silvia Estevez"
$ws.Range("B4").Value = "Departamento Administrativo"

$ws.Range("G2").Select()
